$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B (First Name) for data rows 2-5 to the de-spaced, non-suffixed
# "Nguyen Trong" value, matching the target sharedStrings cleanup.
$ws.Range("B2").Value = "Nguyen Trong"
$ws.Range("B3").Value = "Nguyen Trong"
$ws.Range("B4").Value = "Nguyen Trong"
$ws.Range("B5").Value = "Nguyen Trong"

# Row 5 (previously the "Nhat3" record) now reuses "Nhat" in column A.
$ws.Range("A5").Value = "Nhat"

# The old row 6 ("Nhat4" record) is removed entirely.
$ws.Rows.Item(6).Delete()
